$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at row 147 (pushing the
# existing rows 147-209 down to 148-210, which is why every row below
# shows up as "shifted" in the diff). The new row carries this week's
# Acelga price observation.
$ws.Rows.Item(147).Insert()

$ws.Cells.Item(147, 1).Value = 7
$ws.Cells.Item(147, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(147, 3).Value = "Ñuble"
$ws.Cells.Item(147, 4).Value = 44609
$ws.Cells.Item(147, 5).Value = 16
$ws.Cells.Item(147, 6).Value = 100112009
$ws.Cells.Item(147, 7).Value = "Acelga"
$ws.Cells.Item(147, 8).Value = "Sin especificar"
$ws.Cells.Item(147, 9).Value = "Primera"
$ws.Cells.Item(147, 10).Value = 100
$ws.Cells.Item(147, 11).Value = 400
$ws.Cells.Item(147, 12).Value = 450
$ws.Cells.Item(147, 13).Value = 425
$ws.Cells.Item(147, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(147, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(147, 16).Value = 425
$ws.Cells.Item(147, 17).Value = 1
$ws.Cells.Item(147, 18).Value = "Hortaliza"
